$d = $word.ActiveDocument

$pairs = @(
    @("85-17=", "96-62="),
    @("15+30=", "48-8="),
    @("11+70=", "40+16="),
    @("20+67=", "32-12="),
    @("73-30=", "45+36="),
    @("65-36=", "32-26="),
    @("33+20=", "38+37="),
    @("88-79=", "31+31="),
    @("8+90=", "37+19="),
    @("9-4=", "6+5="),
    @("84+8=", "90+2="),
    @("56+19=", "27+37="),
    @("18+58=", "15+23="),
    @("26+38=", "64-17="),
    @("70-63=", "69+14="),
    @("72+22=", "56-55="),
    @("78-36=", "6+54="),
    @("58-34=", "34+13="),
    @("78-12=", "45+12="),
    @("38-3=", "6+78="),
    @("70-4=", "41+57="),
    @("51-25=", "9+83="),
    @("19+12=", "43+43="),
    @("43-2=", "0+91="),
    @("61-27=", "81-46="),
    @("59+4=", "60-21="),
    @("94-92=", "65+20="),
    @("14+16=", "14+45="),
    @("72-15=", "90-36="),
    @("34-28=", "93-2="),
    @("12+81=", "23+40="),
    @("45-32=", "89-77="),
    @("88+10=", "33+61="),
    @("78-21=", "44-35="),
    @("20+74=", "15+66="),
    @("18-5=", "5+57="),
    @("32+42=", "59-26="),
    @("47+0=", "39-30="),
    @("98-77=", "27+7="),
    @("7+56=", "74-63="),
    @("68-13=", "92-8="),
    @("60-0=", "94-65="),
    @("46+33=", "7-4="),
    @("23+71=", "85-14="),
    @("44-24=", "20+71="),
    @("73-3=", "13-0="),
    @("97-12=", "10+0="),
    @("7+9=", "57-27="),
    @("62-4=", "47+19="),
    @("53+41=", "19+11="),
    @("19+36=", "50-4="),
    @("83-38=", "57-53="),
    @("72-71=", "99-38="),
    @("77-29=", "62-6="),
    @("47+10=", "78-1="),
    @("86-4=", "41-23="),
    @("63-38=", "39+34="),
    @("98+0=", "38-28="),
    @("73+21=", "24+6="),
    @("96-6=", "92-71="),
    @("41+13=", "75-65="),
    @("99-81=", "15-6="),
    @("75-57=", "25+50="),
    @("86-42=", "89-15="),
    @("9-0=", "1+40="),
    @("49+4=", "32+5="),
    @("15+60=", "74-44="),
    @("99-95=", "14-6="),
    @("12+38=", "81+0="),
    @("27+53=", "69-2="),
    @("17+60=", "54-31="),
    @("4-4=", "32+51="),
    @("1+24=", "83-33="),
    @("26+50=", "69+17="),
    @("79-15=", "1+96="),
    @("10+51=", "7+54="),
    @("87-54=", "63-32="),
    @("38-10=", "74-4="),
    @("42-15=", "74+10="),
    @("60-22=", "1+2="),
    @("60-34=", "6+43="),
    @("33+9=", "98-45="),
    @("55-13=", "49-0="),
    @("33+40=", "23+27="),
    @("48+3=", "48-22="),
    @("26-20=", "36+13="),
    @("67-18=", "55-7="),
    @("44-25=", "74-65="),
    @("61+7=", "81-58="),
    @("58+17=", "15-12="),
    @("57+29=", "16+2="),
    @("86-69=", "48+1="),
    @("33+36=", "56-48="),
    @("15+81=", "41-20="),
    @("45+5=", "55-9="),
    @("64+10=", "38+55="),
    @("95-10=", "65-58="),
    @("26+31=", "17+67="),
    @("67-42=", "79-57="),
    @("86-0=", "10+26="),
)

foreach ($p in $pairs) {
    $old = $p[0]
    $new = $p[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"